$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new responsible person ("felelős") for the "Módosítás" task row (row 12).
$ws.Range("C12").Value = "Tiha"

# Move the active selection to D14, matching the saved view state.
$null = $ws.Range("D14").Select()
